$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final data for rows 2..42 (Date serial, Channel, Metric, Value).
# The Channel shared strings "TV" / "Radio" swap roles in this edit, and
# the table grows from 34 data rows (2-35) to 41 data rows (2-42).
$data = @(
  @(45950, "Radio", "Spend", 61),
  @(45901, "TV", "Spend", 57),
  @(45936, "Radio", "Spend", 77),
  @(45950, "Radio", "GRPs", 2),
  @(45929, "Radio", "Spend", 182),
  @(45950, "TV", "Spend", 75),
  @(45915, "Radio", "Spend", 181),
  @(45950, "Radio", "Spend", 61),
  @(45908, "TV", "Spend", 190),
  @(45922, "TV", "Spend", 194),
  @(45915, "TV", "GRPs", 3),
  @(45957, "TV", "GRPs", 5),
  @(45943, "Radio", "GRPs", 5),
  @(45957, "Radio", "GRPs", 3),
  @(45950, "TV", "Spend", 75),
  @(45908, "TV", "GRPs", 9),
  @(45957, "TV", "GRPs", 5),
  @(45936, "TV", "GRPs", 4),
  @(45922, "Radio", "GRPs", 6),
  @(45901, "TV", "Spend", 57),
  @(45915, "TV", "GRPs", 3),
  @(45943, "TV", "Spend", 94),
  @(45908, "TV", "GRPs", 9),
  @(45908, "Radio", "GRPs", 5),
  @(45957, "Radio", "Spend", 63),
  @(45936, "Radio", "Spend", 77),
  @(45901, "Radio", "Spend", 132),
  @(45922, "Radio", "Spend", 178),
  @(45929, "TV", "Spend", 107),
  @(45964, "Radio", "GRPs", 8),
  @(45929, "TV", "Spend", 107),
  @(45964, "Radio", "GRPs", 8),
  @(45908, "Radio", "GRPs", 5),
  @(45915, "Radio", "Spend", 181),
  @(45922, "TV", "GRPs", 6),
  @(45929, "Radio", "GRPs", 3),
  @(45943, "Radio", "Spend", 181),
  @(45943, "TV", "GRPs", 4),
  @(45908, "TV", "Spend", 190),
  @(45922, "TV", "GRPs", 6),
  @(45964, "TV", "Spend", 177)
)

$rowIndex = 2
foreach ($rec in $data) {
  $ws.Cells.Item($rowIndex, 1).Value = $rec[0]
  $ws.Cells.Item($rowIndex, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
  $ws.Cells.Item($rowIndex, 2).Value = $rec[1]
  $ws.Cells.Item($rowIndex, 3).Value = $rec[2]
  $ws.Cells.Item($rowIndex, 4).Value = $rec[3]
  $rowIndex = $rowIndex + 1
}
